$wb = $excel.ActiveWorkbook

# Rename the worksheets (tables -> sheet tabs)
$wb.Worksheets.Item("productsList").Name = "products"
$wb.Worksheets.Item("categoriesList").Name = "categories"

# Make "categories" the active sheet/tab
$wb.Worksheets.Item("categories").Activate()
